$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.139.80"
$ws.Range("E2").Value = "  -1.67%  "

# Row 3
$ws.Range("D3").Value = "1.655.91"
$ws.Range("E3").Value = "  -1.78%  "

# Row 4
$ws.Range("E4").Value = "  +0.48%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.10"
$ws.Range("E5").Value = "  +0.73%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5202"
$ws.Range("E6").Value = "  -2.79%  "

# Row 7
$ws.Range("E7").Value = "  +0.47%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2679"
$ws.Range("E8").Value = "  -0.07%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06326"
$ws.Range("E9").Value = "  -1.61%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.13"
$ws.Range("E10").Value = "  -0.89%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07752"
$ws.Range("E11").Value = "  -0.09%  "

# Row 12
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.436"
$ws.Range("E12").Value = "  -1.39%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.641.84"
$ws.Range("E13").Value = "  -2.62%  "

# Row 14
$ws.Range("D14").Value = "1.883.18"
$ws.Range("E14").Value = "  -1.70%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5475"
$ws.Range("E15").Value = "  -2.77%  "

# Row 16
$ws.Range("E16").Value = "  -1.81%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.86"
$ws.Range("E17").Value = "  -2.04%  "

# Row 18
$ws.Range("D18").Value = "26.194.73"
$ws.Range("E18").Value = "  -1.47%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.005"
$ws.Range("E19").Value = "  +0.53%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.676"
$ws.Range("E20").Value = "  -3.07%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "192.68"
$ws.Range("E21").Value = "  -1.16%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.17"
$ws.Range("E22").Value = "  -2.33%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.112"
$ws.Range("E23").Value = "  -4.52%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.007"
$ws.Range("E24").Value = "  +0.57%  "

# Row 25
$ws.Range("E25").Value = "  -4.49%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1238"
$ws.Range("E26").Value = "  -3.07%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.256"
$ws.Range("E27").Value = "  -3.21%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.14"
$ws.Range("E28").Value = "  -0.72%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.411"
$ws.Range("E29").Value = "  -0.67%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06047"
$ws.Range("E30").Value = "  -1.40%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.288"
$ws.Range("E31").Value = "  +0.71%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.554"
$ws.Range("E32").Value = "  -1.45%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.346"
$ws.Range("E33").Value = "  -3.40%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.654"
$ws.Range("E34").Value = "  -3.05%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9820"
$ws.Range("E35").Value = "  -3.64%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.411"
$ws.Range("E36").Value = "  -0.21%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.776"
$ws.Range("E37").Value = "  -0.85%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5927"
$ws.Range("E38").Value = "  +3.42%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01596"
$ws.Range("E39").Value = "  -3.09%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.963"
$ws.Range("E40").Value = "  +0.12%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8655"
$ws.Range("E41").Value = "  -0.42%  "

# Row 42
$ws.Range("E42").Value = "  +0.25%  "

# Row 43
$ws.Range("D43").Value = "1.035.62"
$ws.Range("E43").Value = "  -1.68%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.91"
$ws.Range("E44").Value = "  -0.40%  "

# Row 45
$ws.Range("D45").Value = "1.798.00"
$ws.Range("E45").Value = "  -2.06%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.24"
$ws.Range("E46").Value = "  +0.02%  "

# Row 47
$ws.Range("E47").Value = "  -0.33%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.006"
$ws.Range("E48").Value = "  +0.59%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.140"
$ws.Range("E49").Value = "  -0.59%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05180"
$ws.Range("E50").Value = "  -0.42%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.476"
$ws.Range("E51").Value = "  +4.01%  "

